# end day commit, modificacion ventana reportes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Servicios")

# Widen column A slightly (stored sheet width 13 -> 14 "characters").
# Excel's ColumnWidth COM property carries a fixed +5/6 char padding versus
# the raw OOXML <col width> value, so back that padding out here.
$ws.Columns.Item(1).ColumnWidth = 14 - (5/6)

# Add the new record row
$ws.Cells.Item(2, 1).Value = 123231232321
$ws.Cells.Item(2, 2).Value = "Canalo"
$ws.Cells.Item(2, 3).Value = 30000
